# "Add files via upload" -- adds a new "Additional Nodes" sheet with a
# name/coordinates lookup table, tweaks a couple of node names/coords on
# the "Coordinates" sheet, and leaves the selection/active-sheet state the
# way the author last left it (on the new sheet).

$wb = $excel.ActiveWorkbook

$wsFinal = $wb.Worksheets.Item(1)        # "Final Nodes"
$wsCoord = $wb.Worksheets.Item(2)        # "Coordinates"

# --- Coordinates sheet: fix a couple of swapped node names/coords -------
$wsCoord.Range("B33").Value = "Governor Santiago"
$wsCoord.Range("C33").Value = "14.692512,120.964504"
$wsCoord.Range("B34").Value = "Dalandanan"
$wsCoord.Range("C34").Value = "14.704206,120.961396"

# Disambiguate the two "Recto" nodes.
$wsCoord.Range("B47").Value = "Recto 1"
$wsCoord.Range("B127").Value = "Recto 2"

# --- Add the new "Additional Nodes" sheet, after "Coordinates" ----------
$wsAdd = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsAdd.Name = "Additional Nodes"

$wsAdd.Columns.Item(2).ColumnWidth = 21.16666667
$wsAdd.Columns.Item(3).ColumnWidth = 37.87760417

$wsAdd.Cells.Item(1, 2).Value = "name"
$wsAdd.Cells.Item(1, 3).Value = "coordinates"

$names = @(
    "Heritage",
    "Baclaran Market",
    "Bayview",
    "City of Dreams",
    "Roosevelt Ave",
    "Road 20",
    "Quirino Mindanao Ave",
    "Kingspoint",
    "Bernardino Ave",
    "VMMC",
    "Belfa St",
    "Yakal",
    "Puregold North Comm",
    "DAR",
    "Sto. Domingo",
    "Vicente Cruz",
    "Manila City Hall",
    "Arellano",
    "Nagtahan",
    "Old Sta. Mesa",
    "Madison",
    "Robinsons Magnolia",
    "Robinsons Galleria",
    "Frontera Drive",
    "11th Ave",
    "Paseo de Roxas",
    "Malugay",
    "Pangilinan Virgina St",
    "Mira Nila Homes",
    "Calle Industria",
    "Tiendesitas"
)

$coords = @(
    "14.536998093374049, 120.9924570648909",
    "14.533005505030609, 120.9927308250385",
    "14.579035647463236, 120.97803135023806",
    "14.523872732491139, 120.99061951083875",
    "14.657982825156735, 121.01969861288812",
    "14.671273811247595, 121.03217831955243",
    "14.690499257703188, 121.02813257092167",
    "14.692769601459757, 121.02981464639207",
    "14.696187685621602, 121.03212261042485",
    "14.659467490994714, 121.03620060237702",
    "14.732783321994267, 121.05563408893393",
    "14.71486779802974, 121.0575950782413",
    "14.704742250000209, 121.08052356792064",
    "14.654327028901207, 121.04954215381416",
    "14.626101156176961, 121.01057906600056",
    "14.61173847794976, 120.99491038842679",
    "14.590618921110261, 120.9804200537169",
    "14.600893621197788, 120.99697330364152",
    "14.60134291847372, 121.00011338322857",
    "14.60249840362917, 121.01198800112263",
    "14.612224302534282, 121.03145456901798",
    "14.615401804408375, 121.03707864746254",
    "14.592286689297485, 121.05938012688362",
    "14.589679338910353, 121.07729608098899",
    "14.560472876863182, 121.05689035696875",
    "14.559872097133809, 121.03057206123688",
    "14.561983401583833, 121.02014405720924",
    "14.670914630311772, 121.03919724171247",
    "14.673155525077044, 121.05896243260233",
    "14.605400688042714, 121.07897682628223",
    "14.585556715852753, 121.07862019578211"
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $wsAdd.Cells.Item($row, 1).Value = $i + 1
    $wsAdd.Cells.Item($row, 2).Value = $names[$i]
    $wsAdd.Cells.Item($row, 3).Value = $coords[$i]
}

# --- Restore per-sheet selection state, visiting sheets in the order ----
# --- the author last touched them, ending on the new sheet. -------------
$wsFinal.Activate()
$wsFinal.Range("H18").Select()

$wsCoord.Activate()
$wsCoord.Range("B47").Select()

$wsAdd.Activate()
$wsAdd.Range("I16").Select()
